$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.708836333333334
$ws.Range("H2").Value = 8.126509
$ws.Range("I2").Value = 0.0171826329450544
$ws.Range("J2").Value = 0.0171826329450544
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03970866666666666
$ws.Range("N2").Value = 0.119126
$ws.Range("Q2").Value = 0.1075642790148889
$ws.Range("R2").Value = 0.968078511134
$ws.Range("S2").Value = 0.0171826329450544
$ws.Range("T2").Value = 0.0171826329450544

# Row 3
$ws.Range("I3").Value = 0.6368977723762839
$ws.Range("J3").Value = 0.6368977723762839
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03970866666666666
$ws.Range("N3").Value = 0.119126
$ws.Range("Q3").Value = 3.987017002045778
$ws.Range("R3").Value = 35.883153018412
$ws.Range("S3").Value = 0.6368977723762839
$ws.Range("T3").Value = 0.6368977723762839

# Row 4
$ws.Range("G4").Value = 54.53410833333334
$ws.Range("H4").Value = 163.602325
$ws.Range("I4").Value = 0.3459195946786617
$ws.Range("J4").Value = 0.3459195946786617
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03970866666666666
$ws.Range("N4").Value = 0.119126
$ws.Range("Q4").Value = 2.165476729772222
$ws.Range("R4").Value = 19.48929056795
$ws.Range("S4").Value = 0.3459195946786617
$ws.Range("T4").Value = 0.3459195946786617
